$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell G4 value from "CURA Healthcare" to "CURA Healthcare Service"
$ws.Range("G4").Value = "CURA Healthcare Service"
